$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update T2 (column C) and T3 (column D) grades for each student row
$ws.Range("C2").Value = 1.25
$ws.Range("D2").Value = 1

$ws.Range("C3").Value = 1.25
$ws.Range("D3").Value = 1.25

$ws.Range("C4").Value = 1.25
$ws.Range("D4").Value = 1.25

$ws.Range("C5").Value = 1.25
$ws.Range("D5").Value = 0

$ws.Range("C6").Value = 1.25
$ws.Range("D6").Value = 1.25

# Update selected cell to D6
$ws.Range("D6").Select()
